$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) cells must stay as text (matching the original inline-string data),
# rather than being auto-converted to numbers by Excel when assigned a
# numeric-looking string. Force text format per-cell before assignment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.108.18'
$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.674.87'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.67'
$ws.Range("E5").Value = '  -2.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5272'
$ws.Range("E6").Value = '  -4.51%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2673'
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06299'
$ws.Range("E9").Value = '  -3.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.28'
$ws.Range("E10").Value = '  -3.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07608'
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.686.62'
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.500'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5685'
$ws.Range("E14").Value = '  -2.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008146'
$ws.Range("E15").Value = '  -3.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.72'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.125.37'
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.833'
$ws.Range("E19").Value = '  -2.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.60'
$ws.Range("E20").Value = '  -3.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '188.79'
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.190'
$ws.Range("E22").Value = '  -0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.006'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.94'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1252'
$ws.Range("E25").Value = '  -5.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.633'
$ws.Range("E26").Value = '  -3.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.80'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06366'
$ws.Range("E28").Value = '  +0.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.355'
$ws.Range("E29").Value = '  -2.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.295'
$ws.Range("E30").Value = '  -2.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.537'
$ws.Range("E31").Value = '  -1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.528'
$ws.Range("E32").Value = '  -1.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.660'
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("E34").Value = '  -2.79%  '
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6022'
$ws.Range("E36").Value = '  -3.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.130'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01622'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.090.95'
$ws.Range("E40").Value = '  -1.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8703'
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("E42").Value = '  -1.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.86'
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.828.24'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05251'
$ws.Range("E48").Value = '  -0.38%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.976'
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.940'
$ws.Range("E51").Value = '  -2.31%  '
